$wb = $excel.ActiveWorkbook

# ---- WELL_GEO (sheet3): add LITH_CODE column (D) ----
$wsGeo = $wb.Worksheets.Item("WELL_GEO")
$wsGeo.Range("D1").Value = "LITH_CODE"
$geoValues = @(100,101,101,102,101,103,101,104,105)
for ($i = 0; $i -lt $geoValues.Length; $i++) {
    $wsGeo.Cells.Item($i + 2, 4).Value = $geoValues[$i]
}

# ---- WELL_ALT (sheet4): rename Lithology -> Alteration, add ALT_CODE column (D) ----
$wsAlt = $wb.Worksheets.Item("WELL_ALT")
$wsAlt.Range("C1").Value = "Alteration"
$wsAlt.Range("D1").Value = "ALT_CODE"
$altValues = @(200,201,200,200,200,200,202,203,200,203)
for ($i = 0; $i -lt $altValues.Length; $i++) {
    $wsAlt.Cells.Item($i + 2, 4).Value = $altValues[$i]
}

# ---- GCHRON (sheet6): update selection (done before activating WELL_ALT so it
#      ends up deselected as the active tab) ----
$wsGchron = $wb.Worksheets.Item("GCHRON")
[void]$wsGchron.Range("S11").Select()

# ---- Restore selections on WELL_GEO / WELL_ALT, then make WELL_ALT the active tab ----
[void]$wsGeo.Range("C1").Select()
[void]$wsAlt.Range("D12").Select()
[void]$wsAlt.Activate()
